$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats / xlPasteValues constants
$xlPasteFormats = -4122
$xlPasteValues = -4163

# -----------------------------------------------------------------------
# 1) "Ativação:" date text: 01/01/2016 -> 01/01/2023
#    Cells B8/C8 (and B15/C15, which reuse the same text) hold this value.
#    Excel's COM layer auto-parses a "dd/mm/yyyy"-looking string into a
#    date serial when assigned through .Value, so we stage the literal
#    text once in an unused scratch cell forced to Text format, copy just
#    the VALUE from there into each target cell (keeping the target's
#    existing number format/style untouched), then drop the scratch cell.
# -----------------------------------------------------------------------
$newDate = "01/01/2023"

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = $newDate

foreach ($addr in @("B8", "C8", "B15", "C15")) {
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteValues)
}
$excel.CutCopyMode = 0
$scratch.Clear()

# -----------------------------------------------------------------------
# 2) New English paragraphs added alongside existing labels.
#    B/C11, B/C14 and B/C16 were previously empty (no <c> at all), so a
#    plain .Value assignment would inherit column A's style instead of
#    the B/C column's wrap-text style; paste the format from a populated
#    sibling row (13) of the same column right after setting the value.
# -----------------------------------------------------------------------
$objectivesText = "To provide the incoming student of Physical Engineering with practical knowledge of electronics and physical computing with Arduino microcontroller aiming its application in scientific and technological projects"
$shortSyllabusText = "Introduction to Arduino. Analog and digital electronics concepts. Assembly of basic electronic circuits. Programming and control of electronic circuits in C language. Application and development of projects based on Arduino."
$syllabusText = "Introduction to the Arduino microcontroller: history, types and resources. Practical workshop: installation and configuration of the Arduino IDE.Basic concepts of electronics: operation of the breadboard, electronic components and instruments, measurements with a multimeter and oscilloscope. Electrical quantities: resistance, voltage and current. Workshop: assembly of electronic circuits.Introduction to the Wiring programming language based on C/C++. Data types, basic syntax, flow control, standard library functions. main librariesArduino Inputs and Outputs. Analog and digital signals.Device control using PWM.Analog electronics. Arduino analog-to-digital converters.Workshop: Reading sensor data. Serial/USB communication with PC. Using the IDE's Serial Monitor.DC motor and servo motor control with PWM. Power control with relay and SSR.Advanced topics: Ethernet communication with Arduino. Wireless communication via Bluetooth.Data storage using ATMega328 EEPROM and SD memory card.Quality software development.Development of projects using Arduino microcontroller."

$newCells = @{
    "B11" = $objectivesText
    "C11" = $objectivesText
    "B14" = $shortSyllabusText
    "C14" = $shortSyllabusText
    "B16" = $syllabusText
    "C16" = $syllabusText
}

foreach ($addr in $newCells.Keys) {
    $col = $addr.Substring(0, 1)
    $ws.Range($addr).Value = $newCells[$addr]
    $ws.Range($col + "13").Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0
